$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row of the sheet (data currently spans A1:AC62)
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# New columns go right after the existing last column (AC -> AD, AE, AF)
$lastCol = $usedRange.Columns.Count
$winsCol = $lastCol + 1
$lossesCol = $lastCol + 2
$tiesCol = $lastCol + 3

# Header row (row 1): add Wins / Losses / Ties headers
$ws.Cells.Item(1, $winsCol).Value = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value = "Ties"

# Copy the existing header style (from the last original header cell) onto the new headers
$lastHeaderCell = $ws.Cells.Item(1, $lastCol)
$lastHeaderCell.Copy()
$newHeaderRange = $ws.Range($ws.Cells.Item(1, $winsCol), $ws.Cells.Item(1, $tiesCol))
$newHeaderRange.PasteSpecial(-4122)  # xlPasteFormats

# Fill the team's win/loss/tie record for every data row (row 2 through the last row)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value = 67
    $ws.Cells.Item($r, $lossesCol).Value = 95
    $ws.Cells.Item($r, $tiesCol).Value = 0
}
